$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted right before the existing row 114,
# pushing all following rows (old 114..201) down by one (new 115..202).
$ws.Rows("114:114").Insert()

$ws.Range("A114").Value = 3
$ws.Range("B114").Value = "Femacal de La Calera"
$ws.Range("C114").Value = "Coquimbo"
$ws.Range("D114").Value = 44447
$ws.Range("E114").Value = 5
$ws.Range("F114").Value = 100112031
$ws.Range("G114").Value = "Poroto verde"
$ws.Range("H114").Value = "Magnum"
$ws.Range("I114").Value = "Primera"
$ws.Range("J114").Value = 73
$ws.Range("K114").Value = 34000
$ws.Range("L114").Value = 35000
$ws.Range("M114").Value = 34521
$ws.Range("N114").Value = "$/malla 25 kilos"
$ws.Range("O114").Value = "Provincia de Quillota"
$ws.Range("P114").Value = 1381
$ws.Range("Q114").Value = 25
$ws.Range("R114").Value = "Hortaliza"
